$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift each scheduled date forward (same month/day, new year) and update the
# French day-of-week label in column B to match the new date.
$ws.Range("A2").Value = 46029
$ws.Range("B2").Value = "mercredi"
$ws.Range("A5").Value = 46036
$ws.Range("B5").Value = "mercredi"
$ws.Range("A8").Value = 46064
$ws.Range("B8").Value = "mercredi"
$ws.Range("A11").Value = 46073
$ws.Range("B11").Value = "vendredi"
$ws.Range("A15").Value = 46090
$ws.Range("B15").Value = "lundi"
$ws.Range("A17").Value = 46091
$ws.Range("B17").Value = "mardi"
$ws.Range("A19").Value = 46092
$ws.Range("B19").Value = "mercredi"
$ws.Range("A22").Value = 46097
$ws.Range("B22").Value = "lundi"
$ws.Range("A24").Value = 46100
$ws.Range("B24").Value = "jeudi"
$ws.Range("A27").Value = 46105
$ws.Range("B27").Value = "mardi"
$ws.Range("A29").Value = 46107
$ws.Range("B29").Value = "jeudi"
$ws.Range("A31").Value = 46108
$ws.Range("B31").Value = "vendredi"
$ws.Range("A34").Value = 46113
$ws.Range("B34").Value = "mercredi"
$ws.Range("A38").Value = 46114
$ws.Range("B38").Value = "jeudi"
$ws.Range("A41").Value = 46119
$ws.Range("B41").Value = "mardi"
$ws.Range("A44").Value = 46120
$ws.Range("B44").Value = "mercredi"
$ws.Range("A47").Value = 46121
$ws.Range("B47").Value = "jeudi"
$ws.Range("A52").Value = 46125
$ws.Range("B52").Value = "lundi"
$ws.Range("A55").Value = 46126
$ws.Range("B55").Value = "mardi"
$ws.Range("A57").Value = 46127
$ws.Range("B57").Value = "mercredi"
$ws.Range("A59").Value = 46128
$ws.Range("B59").Value = "jeudi"
$ws.Range("A62").Value = 46129
$ws.Range("B62").Value = "vendredi"
$ws.Range("A65").Value = 46132
$ws.Range("B65").Value = "lundi"
$ws.Range("A68").Value = 46134
$ws.Range("B68").Value = "mercredi"
$ws.Range("A71").Value = 46135
$ws.Range("B71").Value = "jeudi"
$ws.Range("A74").Value = 46136
$ws.Range("B74").Value = "vendredi"
$ws.Range("A78").Value = 46146
$ws.Range("B78").Value = "lundi"
$ws.Range("A81").Value = 46149
$ws.Range("B81").Value = "jeudi"
$ws.Range("A84").Value = 46154
$ws.Range("B84").Value = "mardi"
$ws.Range("A88").Value = 46160
$ws.Range("B88").Value = "lundi"

# One session's end time moved earlier: 16:45 -> 15:45.
$ws.Range("D30").Value = "15:45"
